$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 46.96651666666666
$ws.Range("H2").Value = 140.89955
$ws.Range("I2").Value = 0.5808027674561179
$ws.Range("J2").Value = 0.5808027674561179
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 161.7750676666667
$ws.Range("N2").Value = 485.325203
$ws.Range("O2").Value = 0.9790864123038654
$ws.Range("P2").Value = 0.9790864123038654
$ws.Range("Q2").Value = 7598.011411817626
$ws.Range("R2").Value = 68382.10270635864
$ws.Range("S2").Value = 0.5686560978447667
$ws.Range("T2").Value = 0.5686560978447667

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 46.96651666666666
$ws.Range("H3").Value = 140.89955
$ws.Range("I3").Value = 0.5808027674561179
$ws.Range("J3").Value = 0.5808027674561179
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.67894
$ws.Range("N3").Value = 2.03682
$ws.Range("O3").Value = 0.004109044356199978
$ws.Range("P3").Value = 0.004109044356199979
$ws.Range("Q3").Value = 31.88744682566666
$ws.Range("R3").Value = 286.987021431
$ws.Range("S3").Value = 0.00238654433368089
$ws.Range("T3").Value = 0.00238654433368089

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 46.96651666666666
$ws.Range("H4").Value = 140.89955
$ws.Range("I4").Value = 0.5808027674561179
$ws.Range("J4").Value = 0.5808027674561179
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.763201333333333
$ws.Range("N4").Value = 5.289604
$ws.Range("O4").Value = 0.01067115280816804
$ws.Range("P4").Value = 0.01067115280816804
$ws.Range("Q4").Value = 82.81142480868887
$ws.Range("R4").Value = 745.3028232781999
$ws.Range("S4").Value = 0.006197835082931123
$ws.Range("T4").Value = 0.006197835082931123

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 46.96651666666666
$ws.Range("H5").Value = 140.89955
$ws.Range("I5").Value = 0.5808027674561179
$ws.Range("J5").Value = 0.5808027674561179
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.013424
$ws.Range("N5").Value = 3.040272
$ws.Range("O5").Value = 0.006133390531766587
$ws.Range("P5").Value = 0.006133390531766588
$ws.Range("Q5").Value = 47.59699518639999
$ws.Range("R5").Value = 428.3729566775999
$ws.Range("S5").Value = 0.003562290194739184
$ws.Range("T5").Value = 0.003562290194739185

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 14.34807866666667
$ws.Range("H6").Value = 43.04423600000001
$ws.Range("I6").Value = 0.1774328689611448
$ws.Range("J6").Value = 0.1774328689611448
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 161.7750676666667
$ws.Range("N6").Value = 485.325203
$ws.Range("O6").Value = 0.9790864123038654
$ws.Range("P6").Value = 0.9790864123038654
$ws.Range("Q6").Value = 2321.161397186657
$ws.Range("R6").Value = 20890.45257467991
$ws.Range("S6").Value = 0.1737221110959491
$ws.Range("T6").Value = 0.1737221110959491

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 14.34807866666667
$ws.Range("H7").Value = 43.04423600000001
$ws.Range("I7").Value = 0.1774328689611448
$ws.Range("J7").Value = 0.1774328689611448
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.67894
$ws.Range("N7").Value = 2.03682
$ws.Range("O7").Value = 0.004109044356199978
$ws.Range("P7").Value = 0.004109044356199979
$ws.Range("Q7").Value = 9.741484529946668
$ws.Range("R7").Value = 87.67336076952002
$ws.Range("S7").Value = 0.0007290795288091624
$ws.Range("T7").Value = 0.0007290795288091624

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.34807866666667
$ws.Range("H8").Value = 43.04423600000001
$ws.Range("I8").Value = 0.1774328689611448
$ws.Range("J8").Value = 0.1774328689611448
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.763201333333333
$ws.Range("N8").Value = 5.289604
$ws.Range("O8").Value = 0.01067115280816804
$ws.Range("P8").Value = 0.01067115280816804
$ws.Range("Q8").Value = 25.29855143583822
$ws.Range("R8").Value = 227.686962922544
$ws.Range("S8").Value = 0.001893413257876032
$ws.Range("T8").Value = 0.001893413257876032

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.34807866666667
$ws.Range("H9").Value = 43.04423600000001
$ws.Range("I9").Value = 0.1774328689611448
$ws.Range("J9").Value = 0.1774328689611448
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.013424
$ws.Range("N9").Value = 3.040272
$ws.Range("O9").Value = 0.006133390531766587
$ws.Range("P9").Value = 0.006133390531766588
$ws.Range("Q9").Value = 14.540687274688
$ws.Range("R9").Value = 130.866185472192
$ws.Range("S9").Value = 0.001088265078510467
$ws.Range("T9").Value = 0.001088265078510467

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.746361333333333
$ws.Range("H10").Value = 5.239084
$ws.Range("I10").Value = 0.02159605538935411
$ws.Range("J10").Value = 0.02159605538935411
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 161.7750676666667
$ws.Range("N10").Value = 485.325203
$ws.Range("O10").Value = 0.9790864123038654
$ws.Range("P10").Value = 0.9790864123038654
$ws.Range("Q10").Value = 282.5177228704503
$ws.Range("R10").Value = 2542.659505834052
$ws.Range("S10").Value = 0.02114440439107827
$ws.Range("T10").Value = 0.02114440439107827

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.746361333333333
$ws.Range("H11").Value = 5.239084
$ws.Range("I11").Value = 0.02159605538935411
$ws.Range("J11").Value = 0.02159605538935411
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.67894
$ws.Range("N11").Value = 2.03682
$ws.Range("O11").Value = 0.004109044356199978
$ws.Range("P11").Value = 0.004109044356199979
$ws.Range("Q11").Value = 1.185674563653333
$ws.Range("R11").Value = 10.67107107288
$ws.Range("S11").Value = 0.00008873914951380763
$ws.Range("T11").Value = 0.00008873914951380764

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.746361333333333
$ws.Range("H12").Value = 5.239084
$ws.Range("I12").Value = 0.02159605538935411
$ws.Range("J12").Value = 0.02159605538935411
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.763201333333333
$ws.Range("N12").Value = 5.289604
$ws.Range("O12").Value = 0.01067115280816804
$ws.Range("P12").Value = 0.01067115280816804
$ws.Range("Q12").Value = 3.079186631415111
$ws.Range("R12").Value = 27.712679682736
$ws.Range("S12").Value = 0.0002304548071134587
$ws.Range("T12").Value = 0.0002304548071134587

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.746361333333333
$ws.Range("H13").Value = 5.239084
$ws.Range("I13").Value = 0.02159605538935411
$ws.Range("J13").Value = 0.02159605538935411
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.013424
$ws.Range("N13").Value = 3.040272
$ws.Range("O13").Value = 0.006133390531766587
$ws.Range("P13").Value = 0.006133390531766588
$ws.Range("Q13").Value = 1.769804487872
$ws.Range("R13").Value = 15.928240390848
$ws.Range("S13").Value = 0.0001324570416485713
$ws.Range("T13").Value = 0.0001324570416485713

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 17.80387266666667
$ws.Range("H14").Value = 53.411618
$ws.Range("I14").Value = 0.2201683081933832
$ws.Range("J14").Value = 0.2201683081933832
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 161.7750676666667
$ws.Range("N14").Value = 485.325203
$ws.Range("O14").Value = 0.9790864123038654
$ws.Range("P14").Value = 0.9790864123038654
$ws.Range("Q14").Value = 2880.222705378717
$ws.Range("R14").Value = 25922.00434840845
$ws.Range("S14").Value = 0.2155637989720713
$ws.Range("T14").Value = 0.2155637989720713

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 17.80387266666667
$ws.Range("H15").Value = 53.411618
$ws.Range("I15").Value = 0.2201683081933832
$ws.Range("J15").Value = 0.2201683081933832
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.67894
$ws.Range("N15").Value = 2.03682
$ws.Range("O15").Value = 0.004109044356199978
$ws.Range("P15").Value = 0.004109044356199979
$ws.Range("Q15").Value = 12.08776130830667
$ws.Range("R15").Value = 108.78985177476
$ws.Range("S15").Value = 0.0009046813441961188
$ws.Range("T15").Value = 0.0009046813441961189

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 17.80387266666667
$ws.Range("H16").Value = 53.411618
$ws.Range("I16").Value = 0.2201683081933832
$ws.Range("J16").Value = 0.2201683081933832
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.763201333333333
$ws.Range("N16").Value = 5.289604
$ws.Range("O16").Value = 0.01067115280816804
$ws.Range("P16").Value = 0.01067115280816804
$ws.Range("Q16").Value = 31.39181202436355
$ws.Range("R16").Value = 282.526308219272
$ws.Range("S16").Value = 0.002349449660247428
$ws.Range("T16").Value = 0.002349449660247428

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 17.80387266666667
$ws.Range("H17").Value = 53.411618
$ws.Range("I17").Value = 0.2201683081933832
$ws.Range("J17").Value = 0.2201683081933832
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.013424
$ws.Range("N17").Value = 3.040272
$ws.Range("O17").Value = 0.006133390531766587
$ws.Range("P17").Value = 0.006133390531766588
$ws.Range("Q17").Value = 18.042871853344
$ws.Range("R17").Value = 162.385846680096
$ws.Range("S17").Value = 0.001350378216868365
$ws.Range("T17").Value = 0.001350378216868365
